# Recalculated market-board profit figures for several Leve rows
# across multiple job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Columns H..N hold price/profit figures recomputed by the scheduled
# market-data runner; this script pokes the refreshed values in place.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 530.3333
$ws.Range("I12").Value = 778
$ws.Range("K12").Value = 778
$ws.Range("M12").Value = -608
# row 43
$ws.Range("H43").Value = 1100
$ws.Range("I43").Value = 1100
$ws.Range("J43").Value = 1100
$ws.Range("K43").Value = 1100
$ws.Range("L43").Value = 1100
$ws.Range("M43").Value = -1031
$ws.Range("N43").Value = -1238
# row 48
$ws.Range("H48").Value = 1875
$ws.Range("J48").Value = 1999.6666
$ws.Range("L48").Value = 5998.9998
$ws.Range("N48").Value = -6582.9998
# row 56
$ws.Range("H56").Value = 1875
$ws.Range("J56").Value = 1999.6666
$ws.Range("L56").Value = 5998.9998
$ws.Range("N56").Value = -7066.9998
# row 58
$ws.Range("H58").Value = 381.66666
$ws.Range("J58").Value = 396.66666
$ws.Range("L58").Value = 1189.99998
$ws.Range("N58").Value = -1489.99998
# row 100
$ws.Range("H100").Value = 1274.7142
$ws.Range("I100").Value = 896.3
$ws.Range("J100").Value = 2220.75
$ws.Range("K100").Value = 896.3
$ws.Range("L100").Value = 2220.75
$ws.Range("M100").Value = -355.3
$ws.Range("N100").Value = -3302.75
# row 138
$ws.Range("H138").Value = 5086.436
$ws.Range("I138").Value = 2481.5715
$ws.Range("K138").Value = 7444.7145
$ws.Range("M138").Value = -2304.7145
# row 141
$ws.Range("H141").Value = 2530.3333
$ws.Range("I141").Value = 2530.3333
$ws.Range("K141").Value = 7590.999899999999
$ws.Range("M141").Value = -2410.999899999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value = 919.8
$ws.Range("I4").Value = 550
$ws.Range("J4").Value = 1166.3334
$ws.Range("K4").Value = 550
$ws.Range("L4").Value = 1166.3334
$ws.Range("M4").Value = -434
$ws.Range("N4").Value = -1398.3334
# row 5
$ws.Range("H5").Value = 686.5
$ws.Range("I5").Value = 686.5
$ws.Range("K5").Value = 686.5
$ws.Range("M5").Value = -574.5
# row 6
$ws.Range("H6").Value = 4599.3335
$ws.Range("I6").Value = 5697.5
$ws.Range("J6").Value = 2403
$ws.Range("K6").Value = 5697.5
$ws.Range("L6").Value = 2403
$ws.Range("M6").Value = -5524.5
$ws.Range("N6").Value = -2749
# row 45
$ws.Range("H45").Value = 1320.3334
$ws.Range("I45").Value = 1084.4
$ws.Range("K45").Value = 1084.4
$ws.Range("M45").Value = -707.4000000000001
# row 105
$ws.Range("H105").Value = 200000
$ws.Range("J105").Value = 200000
$ws.Range("L105").Value = 200000
$ws.Range("N105").Value = -206988
# row 122
$ws.Range("H122").Value = 1070.5714
$ws.Range("I122").Value = 798.8
$ws.Range("K122").Value = 2396.4
$ws.Range("M122").Value = 53.60000000000036

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 686.5
$ws.Range("I4").Value = 686.5
$ws.Range("K4").Value = 686.5
$ws.Range("M4").Value = -571.5
# row 94
$ws.Range("H94").Value = 580.25
$ws.Range("I94").Value = 544.8
$ws.Range("K94").Value = 544.8
$ws.Range("M94").Value = -93.79999999999995

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 83.59999999999999
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 83.59999999999999
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 29.40000000000001
$ws.Range("N7").Value = -276
# row 58
$ws.Range("H58").Value = 2264.4
$ws.Range("I58").Value = 1997.25
$ws.Range("J58").Value = 3333
$ws.Range("K58").Value = 1997.25
$ws.Range("L58").Value = 3333
$ws.Range("M58").Value = -1794.25
$ws.Range("N58").Value = -3739
# row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# row 99
$ws.Range("H99").Value = 2062
$ws.Range("J99").Value = 2232.5
$ws.Range("L99").Value = 2232.5
$ws.Range("N99").Value = -5228.5
# row 126
$ws.Range("H126").Value = 2062
$ws.Range("J126").Value = 2232.5
$ws.Range("L126").Value = 6697.5
$ws.Range("N126").Value = -11637.5
# row 136
$ws.Range("H136").Value = 2264.4
$ws.Range("I136").Value = 1997.25
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 5991.75
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -3441.75
$ws.Range("N136").Value = -15099

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 9
$ws.Range("H9").Value = 1718
$ws.Range("I9").Value = 8000
$ws.Range("J9").Value = 147.5
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 442.5
$ws.Range("M9").Value = -23776
$ws.Range("N9").Value = -890.5
# row 12
$ws.Range("H12").Value = 124.25
$ws.Range("J12").Value = 67.333336
$ws.Range("L12").Value = 202.000008
$ws.Range("N12").Value = -548.000008
# row 46
$ws.Range("H46").Value = 1221.4546
$ws.Range("I46").Value = 973.8
$ws.Range("K46").Value = 2921.4
$ws.Range("M46").Value = -2830.4
# row 112
$ws.Range("H112").Value = 9872.111000000001
$ws.Range("I112").Value = 1769.8
$ws.Range("K112").Value = 5309.4
$ws.Range("M112").Value = -4201.4
# row 128
$ws.Range("H128").Value = 557092
$ws.Range("I128").Value = 557092
$ws.Range("K128").Value = 1671276
$ws.Range("M128").Value = -1666296
# row 131
$ws.Range("H131").Value = 2720.75
$ws.Range("I131").Value = 2200
$ws.Range("J131").Value = 2768.0908
$ws.Range("K131").Value = 6600
$ws.Range("L131").Value = 8304.2724
$ws.Range("M131").Value = -1560
$ws.Range("N131").Value = -18384.2724

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 5068.3335
$ws.Range("I80").Value = 1942.5
$ws.Range("K80").Value = 1942.5
$ws.Range("M80").Value = -944.5
# row 83
$ws.Range("H83").Value = 5068.3335
$ws.Range("I83").Value = 1942.5
$ws.Range("K83").Value = 9712.5
$ws.Range("M83").Value = -4720.5
# row 126
$ws.Range("H126").Value = 1999.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 3433.2
$ws.Range("J46").Value = 5500
$ws.Range("L46").Value = 5500
$ws.Range("N46").Value = -5876

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 2997
$ws.Range("M113").Value = -827
